$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) to match the new title
$ws.Name = "Update Group"

# The template now only has a single "Item Group Name" column.
# Put that label in A1 (replacing the old "Item Category Name")...
$ws.Range("A1").Value = "Item Group Name"

# ...and remove the old second column entirely.
$ws.Columns.Item(2).Delete()

# Clear selection back to A1 so no stray selection remains
$ws.Range("A1").Select()
